$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Major") holds the category for each response row, but it is
# only populated on some of the rows. Go through each row and fill in its
# category value, copying the formatting already used elsewhere in the
# response table (e.g. B2) onto the cells that are still blank.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null

$ws.Range("E2").Value  = "Neuroscience"
$ws.Range("E3").Value  = "Neuroscience"
$ws.Range("E4").Value  = "Neuroscience"
$ws.Range("E5").Value  = "History"
$ws.Range("E6").Value  = "History"
$ws.Range("E7").Value  = "History"
$ws.Range("E10").Value = "Neuroscience"

# Update the active selection to reflect where we ended up after going
# through the category column.
$ws.Range("D10").Select()
